$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default/unstyled) used to restore cell style after forcing text format
$plainStyle = $ws.Range('F1').Style

$ws.Range('D2').Value = '37.121.24'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').Value = '2.013.15'
$ws.Range('E3').Value = '  -1.65%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('D4').Style = $plainStyle
$ws.Range('E4').Value = '  +0.75%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.21'
$ws.Range('D5').Style = $plainStyle
$ws.Range('E5').Value = '  -1.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.607'
$ws.Range('D6').Style = $plainStyle
$ws.Range('E6').Value = '  -1.17%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '55.28'
$ws.Range('D8').Style = $plainStyle
$ws.Range('E8').Value = '  -2.65%  '
$ws.Range('E9').Value = '  -2.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0779'
$ws.Range('D10').Style = $plainStyle
$ws.Range('E10').Value = '  -4.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.102'
$ws.Range('D11').Style = $plainStyle
$ws.Range('E11').Value = '  -4.70%  '
$ws.Range('D12').Value = '2.311.19'
$ws.Range('E12').Value = '  -1.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.06'
$ws.Range('D13').Style = $plainStyle
$ws.Range('E13').Value = '  -3.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '19.84'
$ws.Range('D14').Style = $plainStyle
$ws.Range('E14').Value = '  -4.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.21'
$ws.Range('D15').Style = $plainStyle
$ws.Range('E15').Value = '  -1.48%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.738'
$ws.Range('D16').Style = $plainStyle
$ws.Range('E16').Value = '  -2.45%  '
$ws.Range('D17').Value = '2.014.28'
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('D18').Value = '37.041.00'
$ws.Range('E18').Value = '  -0.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.22'
$ws.Range('D19').Style = $plainStyle
$ws.Range('E19').Value = '  +3.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.39'
$ws.Range('D20').Style = $plainStyle
$ws.Range('E20').Value = '  -2.11%  '
$ws.Range('D21').Value = '0.0₃0814'
$ws.Range('E21').Value = '  -3.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '223.80'
$ws.Range('D22').Style = $plainStyle
$ws.Range('E22').Value = '  -1.15%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('E24').Value = '  +2.29%  '
$ws.Range('E25').Value = '  -5.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.60'
$ws.Range('D26').Style = $plainStyle
$ws.Range('E26').Value = '  -1.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.95'
$ws.Range('D27').Style = $plainStyle
$ws.Range('E27').Value = '  -5.92%  '
$ws.Range('E28').Value = '  -3.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.62'
$ws.Range('D29').Style = $plainStyle
$ws.Range('E29').Value = '  -1.89%  '
$ws.Range('E30').Value = '  -7.87%  '
$ws.Range('E31').Value = '  -1.41%  '
$ws.Range('E32').Value = '  -1.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0600'
$ws.Range('D33').Style = $plainStyle
$ws.Range('E33').Value = '  -2.15%  '
$ws.Range('E34').Value = '  -2.57%  '
$ws.Range('E35').Value = '  -3.37%  '
$ws.Range('E36').Value = '  +2.04%  '
$ws.Range('E37').Value = '  +0.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.12'
$ws.Range('D38').Style = $plainStyle
$ws.Range('E38').Value = '  -3.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.32'
$ws.Range('D39').Style = $plainStyle
$ws.Range('E39').Value = '  -1.28%  '
$ws.Range('D40').Value = '1.456.18'
$ws.Range('E40').Value = '  -1.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0212'
$ws.Range('D41').Style = $plainStyle
$ws.Range('E41').Value = '  -3.94%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '94.80'
$ws.Range('D42').Style = $plainStyle
$ws.Range('E42').Value = '  -1.19%  '
$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.76'
$ws.Range('D43').Style = $plainStyle
$ws.Range('E43').Value = '  -4.76%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0905'
$ws.Range('D44').Style = $plainStyle
$ws.Range('E44').Value = '  -4.39%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.07'
$ws.Range('D45').Style = $plainStyle
$ws.Range('E45').Value = '  -5.42%  '
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.12'
$ws.Range('D46').Style = $plainStyle
$ws.Range('E46').Value = '  -2.66%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.14'
$ws.Range('D47').Style = $plainStyle
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.999'
$ws.Range('D48').Style = $plainStyle
$ws.Range('E48').Value = '  -1.99%  '
$ws.Range('E49').Value = '  -0.35%  '
$ws.Range('B50').Value = 'FTXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.76'
$ws.Range('D50').Style = $plainStyle
$ws.Range('E50').Value = '  +1.63%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.198.52'
$ws.Range('E51').Value = '  -1.68%  '
